$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.664.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.962.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.95"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.82"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2955"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06797"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "110.73"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.39"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.964.63"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07743"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.487"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.6924"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "294.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "30.684.92"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.34"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.668"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007704"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.214.13"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.645"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.882"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.93"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.18"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.200"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1077"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.441"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.685"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +16.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.444"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05103"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7794"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.176"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02060"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.737"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.717"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.070"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.32"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.091"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4470"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8754"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.19"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.421"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1279"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.360"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.93"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.83"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.01%  "
